# Daily attendance processing - 2025-11-22 23:21:29
# For every row's "Recorded By" cell (column G) that holds a comma-separated
# list of more than one recorder, reverse the order of the names/emails in
# the list - unless the list includes "admin@admin.com", which is left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "

    if ($parts.Count -le 1) { continue }
    if ($val -like "*admin@admin.com*") { continue }

    $reversed = @()
    for ($i = $parts.Count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $cell.Value = $reversed -join ", "
}
